$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the "L" column (day-name label next to the "APOYO"/capataz column)
# for rows 5-18 with the Spanish number words for days 15-28. These cells
# were previously blank (rows 5,7,8,13,15,16,18) or held leftover test
# strings (rows 6,9,10,11,12,14,17) that get replaced here.
$ws.Range("L5").Value  = "QUINCE"
$ws.Range("L6").Value  = "DIECISEIS"
$ws.Range("L7").Value  = "DIECISIETE"
$ws.Range("L8").Value  = "DIECIOCHO"
$ws.Range("L9").Value  = "DIECINUEVE"
$ws.Range("L10").Value = "VEINTE"
$ws.Range("L11").Value = "VEINTIUNO"
$ws.Range("L12").Value = "VEINTIDOS"
$ws.Range("L13").Value = "VIENTITRES"
$ws.Range("L14").Value = "VEINTICUATRO"
$ws.Range("L15").Value = "VEINTICINCO"
$ws.Range("L16").Value = "VEINTISEIS"
$ws.Range("L17").Value = "VEINTISIETE"
$ws.Range("L18").Value = "VEINTIOCHO"

# Re-enter the TOTALES row formulas as one pass so they collapse back into a
# single shared-formula group (matches the saved workbook behaviour).
$ws.Range("C19:I19").Formula = "=SUM(C5:C18)"

# Move the active selection to L18 (matches the last cell touched).
$ws.Range("L18").Select()
